$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6467574834823608
$ws.Range("B1").Value = 1.009629368782043
$ws.Range("C1").Value = 2.597643136978149
$ws.Range("D1").Value = 6.258825778961182
$ws.Range("E1").Value = 2.12835955619812
